# Trading update: 2026-02-18 10:21:56
$wb = $excel.ActiveWorkbook

# --- Update "All Trades" sheet ---
$trades = $wb.Worksheets.Item("All Trades")

# Row 11 (Trade #10) used to be the "latest open trade" with extra computed
# fields filled in (Capital After, slippages, confidence, entry reason,
# duration). Now that newer trades exist, it becomes a plain open-trade row.
$trades.Range("G11").Value = 0
$trades.Range("K11").Value = ""
$trades.Range("L11").Value = ""
$trades.Range("M11").Value = ""
$trades.Range("N11").Value = ""
$trades.Range("O11").Value = ""
$trades.Range("Q11").Value = ""

# New trade #11 -> row 12
$trades.Range("A12").Value = 11
$trades.Range("B12").Value = "'2026-02-18"
$trades.Range("C12").Value = "10:21:16"
$trades.Range("D12").Value = "MarketMaking"
$trades.Range("E12").Value = "UP"
$trades.Range("F12").Value = 0.55
$trades.Range("G12").Value = 0
$trades.Range("H12").Value = "OPEN"
$trades.Range("I12").Value = 0
$trades.Range("J12").Value = 0
$trades.Range("K12").Value = ""
$trades.Range("L12").Value = ""
$trades.Range("M12").Value = ""
$trades.Range("N12").Value = ""
$trades.Range("O12").Value = ""
$trades.Range("P12").Value = ""
$trades.Range("Q12").Value = ""

# New trade #12 -> row 13
$trades.Range("A13").Value = 12
$trades.Range("B13").Value = "'2026-02-18"
$trades.Range("C13").Value = "10:21:22"
$trades.Range("D13").Value = "MarketMaking"
$trades.Range("E13").Value = "UP"
$trades.Range("F13").Value = 0.58
$trades.Range("G13").Value = 0
$trades.Range("H13").Value = "OPEN"
$trades.Range("I13").Value = 0
$trades.Range("J13").Value = 0
$trades.Range("K13").Value = ""
$trades.Range("L13").Value = ""
$trades.Range("M13").Value = ""
$trades.Range("N13").Value = ""
$trades.Range("O13").Value = ""
$trades.Range("P13").Value = ""
$trades.Range("Q13").Value = ""

# New trade #13 -> row 14 (newest open trade, gets the extra computed fields)
$trades.Range("A14").Value = 13
$trades.Range("B14").Value = "'2026-02-18"
$trades.Range("C14").Value = "10:21:54"
$trades.Range("D14").Value = "MarketMaking"
$trades.Range("E14").Value = "DOWN"
$trades.Range("F14").Value = 0.27
$trades.Range("G14").Value = ""
$trades.Range("H14").Value = "OPEN"
$trades.Range("I14").Value = 0
$trades.Range("J14").Value = 0
$trades.Range("K14").Value = 100
$trades.Range("L14").Value = 0
$trades.Range("M14").Value = 0
$trades.Range("N14").Value = 0.6
$trades.Range("O14").Value = "Normal spread capture: 202 bps"
$trades.Range("P14").Value = ""
$trades.Range("Q14").Value = 0

# --- Update "MarketMaking" sheet (mirrors latest open trade for this strategy) ---
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("A2").Value = 13
$mm.Range("C2").Value = "10:21:54"
$mm.Range("F2").Value = 0.27
$mm.Range("O2").Value = "Normal spread capture: 202 bps"
